# Restrict caption files to individual masterfiles
#
# The manifest's "caption" column group (Caption File / Caption Label /
# Caption Language, columns P:R) and the "structural metadata" column
# group (File / Offset / Label / Date Digitized, columns S:V) swap places:
# the File/Offset/Label/Date Digitized block now comes first (P:S),
# immediately after "Note Type", followed by the Caption File/Label/
# Language block (T:V). Everything from column W onward is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the current ("before") values + number formats for the
#     two column blocks that are swapping, rows 1-5 -------------------
$capturedValues = @{}
$capturedFormats = @{}
foreach ($row in 1..5) {
    foreach ($col in @("P","Q","R","S","T","U","V")) {
        $ref = "$col$row"
        $cell = $ws.Range($ref)
        $capturedValues[$ref] = $cell.Value()
        $capturedFormats[$ref] = $cell.NumberFormat()
    }
}
# the lone stray formatted-but-empty cell at T1 moves to Q1
$t1Format = $ws.Range("T1").NumberFormat()

# --- clear the whole P1:V5 block before rewriting it -----------------
foreach ($row in 1..5) {
    foreach ($col in @("P","Q","R","S","T","U","V")) {
        $ref = "$col$row"
        $cell = $ws.Range($ref)
        $cell.Value = ""
        $cell.NumberFormat = "General"
    }
}

# --- write the old S:V (File/Offset/Label/Date Digitized) block into
#     the new P:S location, row by row -------------------------------
$srcBlock1 = @("S","T","U","V")
$dstBlock1 = @("P","Q","R","S")
foreach ($row in 1..5) {
    for ($i = 0; $i -lt 4; $i++) {
        $srcRef = "$($srcBlock1[$i])$row"
        $dstRef = "$($dstBlock1[$i])$row"
        $val = $capturedValues[$srcRef]
        if ($val -ne $null -and $val -ne "") {
            $ws.Range($dstRef).Value = $val
            $ws.Range($dstRef).NumberFormat = $capturedFormats[$srcRef]
        }
    }
}

# --- write the old P:R (Caption File/Label/Language) block into the
#     new T:V location, row by row ------------------------------------
$srcBlock2 = @("P","Q","R")
$dstBlock2 = @("T","U","V")
foreach ($row in 1..5) {
    for ($i = 0; $i -lt 3; $i++) {
        $srcRef = "$($srcBlock2[$i])$row"
        $dstRef = "$($dstBlock2[$i])$row"
        $val = $capturedValues[$srcRef]
        if ($val -ne $null -and $val -ne "") {
            $ws.Range($dstRef).Value = $val
            $ws.Range($dstRef).NumberFormat = $capturedFormats[$srcRef]
        }
    }
}

# --- the stray formatted empty cell: was T1 (2nd col of old S:V
#     block), now sits at Q1 (2nd col of new P:S block) ---------------
$ws.Range("Q1").NumberFormat = $t1Format

# --- sheet view: active cell / scroll position moved from L3 / D1 to
#     P1 / G1 -----------------------------------------------------------
$ws.Range("P1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 7
